$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 2.745448847643053
$ws.Range("D2").Value = 9.493688876567631

# Row 3
$ws.Range("B3").Value = 6.934574674567631
$ws.Range("C3").Value = -16.12831465753368
$ws.Range("D3").Value = -2.771911005877533

# Row 4
$ws.Range("B4").Value = 2.666952881567631
$ws.Range("C4").Value = -2.27083943853368
$ws.Range("D4").Value = -0.3450489484686868

# Row 5
$ws.Range("B5").Value = -2.75138498353368
$ws.Range("C5").Value = -2.111703357971903
$ws.Range("D5").Value = 2.040358025009379

# Row 6
$ws.Range("B6").Value = -0.9167407425839711
$ws.Range("C6").Value = 1.912948101009379
$ws.Range("D6").Value = -1.764067355060938

# Row 7
$ws.Range("B7").Value = 1.349025510009379
$ws.Range("C7").Value = -1.676605839260361
$ws.Range("D7").Value = 0.9272117483245859

# Row 8
$ws.Range("B8").Value = -1.391220893296087
$ws.Range("C8").Value = 0.2611223103245859
$ws.Range("D8").Value = -0.5508352566021527

# Row 9
$ws.Range("B9").Value = -0.137926112675414
$ws.Range("C9").Value = -0.7839136486021527
$ws.Range("D9").Value = -0.02328789311612689

# Row 10
$ws.Range("B10").Value = -0.3849111476021527
$ws.Range("C10").Value = 1.598289643883873
$ws.Range("D10").Value = 0.1507850357460907

# Row 11
$ws.Range("B11").Value = 0.09167924288387319
$ws.Range("C11").Value = 0.7764434987460908
$ws.Range("D11").Value = -0.6853637869931832

# Row 12
$ws.Range("B12").Value = 0.3758430847460907
$ws.Range("C12").Value = -0.5715147199931834
$ws.Range("D12").Value = 0.7147782794088806

# Row 13
$ws.Range("B13").Value = -0.8121133989931834
$ws.Range("C13").Value = 0.6007027414088806
$ws.Range("D13").Value = 1.071245387738463

# Row 14
$ws.Range("B14").Value = 0.5217354724088806
$ws.Range("C14").Value = 1.012535832738463
$ws.Range("D14").Value = 0.2265041937008457

# Row 15
$ws.Range("B15").Value = 1.200271971738463
$ws.Range("C15").Value = -0.2426956912991542
$ws.Range("D15").Value = -0.716353698

# Row 16
$ws.Range("B16").Value = -0.2129870522991542
$ws.Range("C16").Value = -0.474503149
$ws.Range("D16").Value = 0.1529377914640456

# Row 17
$ws.Range("B17").Value = -0.383419286
$ws.Range("C17").Value = -0.01429268463133643
$ws.Range("D17").Value = 0.3607736698111976

# Row 18
$ws.Range("B18").Value = -0.0881763245104944
$ws.Range("C18").Value = -0.1315562061888024
$ws.Range("D18").Value = 0.4244005293199388

# Row 19
$ws.Range("B19").Value = -0.3266301831888024
$ws.Range("C19").Value = 0.6331622283199388
$ws.Range("D19").Value = 0.4436511628968201

# Row 20
$ws.Range("B20").Value = 0.08648729331993882
$ws.Range("C20").Value = -0.5784485191031798
$ws.Range("D20").Value = 0.07992519753787367

# Row 21
$ws.Range("B21").Value = -0.6655694571031798
$ws.Range("C21").Value = -0.0467553644621263
$ws.Range("D21").Value = -0.2524569298566441

# Row 22 - gains a new D22 cell
$ws.Range("B22").Value = 0.09619904853787367
$ws.Range("C22").Value = 0.02585042314335589
$ws.Range("D22").Value = 0.7115302101128926

# Row 23 - gains a new C23 cell
$ws.Range("B23").Value = 0.2427972171433558
$ws.Range("C23").Value = 0.3054124296933831
